$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 220, shifting existing rows 220-228 down to 223-231
$ws.Rows("220:222").Insert()

# Populate the 3 new rows (220-222) with the new weekly price group (date 44461 / 2021-09-22)
$ws.Cells.Item(220, 1).Value = 2
$ws.Cells.Item(220, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(220, 3).Value = "Coquimbo"
$ws.Cells.Item(220, 4).Value = 44461
$ws.Cells.Item(220, 5).Value = 4
$ws.Cells.Item(220, 6).Value = "Fruta"
$ws.Cells.Item(220, 7).Value = 100101
$ws.Cells.Item(220, 8).Value = "Berries"
$ws.Cells.Item(220, 9).Value = 100112025
$ws.Cells.Item(220, 10).Value = "Frutilla"
$ws.Cells.Item(220, 11).Value = "Sin especificar"
$ws.Cells.Item(220, 12).Value = "Especial"
$ws.Cells.Item(220, 13).Value = 160
$ws.Cells.Item(220, 14).Value = 26000
$ws.Cells.Item(220, 15).Value = 27000
$ws.Cells.Item(220, 16).Value = 26500
$ws.Cells.Item(220, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(220, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(220, 19).Value = 3786
$ws.Cells.Item(220, 20).Value = 7

$ws.Cells.Item(221, 1).Value = 2
$ws.Cells.Item(221, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(221, 3).Value = "Coquimbo"
$ws.Cells.Item(221, 4).Value = 44461
$ws.Cells.Item(221, 5).Value = 4
$ws.Cells.Item(221, 6).Value = "Fruta"
$ws.Cells.Item(221, 7).Value = 100101
$ws.Cells.Item(221, 8).Value = "Berries"
$ws.Cells.Item(221, 9).Value = 100112025
$ws.Cells.Item(221, 10).Value = "Frutilla"
$ws.Cells.Item(221, 11).Value = "Sin especificar"
$ws.Cells.Item(221, 12).Value = "Primera"
$ws.Cells.Item(221, 13).Value = 240
$ws.Cells.Item(221, 14).Value = 23000
$ws.Cells.Item(221, 15).Value = 24000
$ws.Cells.Item(221, 16).Value = 23500
$ws.Cells.Item(221, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(221, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(221, 19).Value = 3357
$ws.Cells.Item(221, 20).Value = 7

$ws.Cells.Item(222, 1).Value = 2
$ws.Cells.Item(222, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(222, 3).Value = "Coquimbo"
$ws.Cells.Item(222, 4).Value = 44461
$ws.Cells.Item(222, 5).Value = 4
$ws.Cells.Item(222, 6).Value = "Fruta"
$ws.Cells.Item(222, 7).Value = 100101
$ws.Cells.Item(222, 8).Value = "Berries"
$ws.Cells.Item(222, 9).Value = 100112025
$ws.Cells.Item(222, 10).Value = "Frutilla"
$ws.Cells.Item(222, 11).Value = "Sin especificar"
$ws.Cells.Item(222, 12).Value = "Segunda"
$ws.Cells.Item(222, 13).Value = 240
$ws.Cells.Item(222, 14).Value = 19000
$ws.Cells.Item(222, 15).Value = 20000
$ws.Cells.Item(222, 16).Value = 19500
$ws.Cells.Item(222, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(222, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(222, 19).Value = 2786
$ws.Cells.Item(222, 20).Value = 7
